$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the FFR_A / FFR_LF columns entirely (columns B:C), shifting
# C_A / A_C left into B:C so the table becomes a 2-variable (3x3) block.
$ws.Range("B1:C3").Delete(-4159)

# Refresh the two remaining data values to match the recomputed results.
$ws.Range("B2").Value = 0.09733919936953835
$ws.Range("C2").Value = 9.720664130955248

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
